# acpt: add scenarios for row height and height_rule
#
# Sets explicit row-height / height-rule (w:trHeight) on the first row of
# four of the tables in the document:
#   Table 2 (1-col, 1440 dxa)        -> auto,    0 dxa
#   Table 3 (1-col, 2880 dxa)        -> atLeast, 2880 dxa
#   Table 4 (3-col, 1440 dxa cols)   -> exact,   4320 dxa (first row only)
#   Table 5 (3-col, 1440 dxa cols)   -> exact,   5760 dxa (first row only)
#
# Row.Height is expressed in points; w:trHeight/@w:val is in twips (dxa),
# i.e. 20 twips per point, so the dxa target is divided by 20 below.
#
# WdRowHeightRule: wdRowHeightAuto = 0, wdRowHeightAtLeast = 1,
#                  wdRowHeightExactly = 2

$d = $word.ActiveDocument

$row1 = $d.Tables.Item(2).Rows.Item(1)
$row1.HeightRule = 0
$row1.Height = 0

$row2 = $d.Tables.Item(3).Rows.Item(1)
$row2.HeightRule = 1
$row2.Height = 2880 / 20

$row3 = $d.Tables.Item(4).Rows.Item(1)
$row3.HeightRule = 2
$row3.Height = 4320 / 20

$row4 = $d.Tables.Item(5).Rows.Item(1)
$row4.HeightRule = 2
$row4.Height = 5760 / 20
